$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (K = strikeouts) values per regenerated save_data
$kValues = @{
    2 = 0
    4 = 1
    5 = 0
    6 = 0
    7 = 0
    8 = 1
    9 = 0
    10 = 0
    11 = 0
    12 = 2
    13 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 2
    26 = 1
    27 = 0
    28 = 3
    29 = 2
    30 = 0
    31 = 2
    34 = 1
    35 = 1
    36 = 0
    37 = 2
    38 = 1
    39 = 1
    41 = 0
    42 = 0
    43 = 1
    45 = 2
    46 = 0
    47 = 1
    48 = 0
    49 = 1
    50 = 2
    51 = 2
    52 = 2
    53 = 0
    54 = 1
    55 = 0
    56 = 1
    57 = 1
    58 = 2
    59 = 2
    60 = 0
    61 = 0
    62 = 0
    63 = 1
    64 = 1
    65 = 2
    66 = 1
    67 = 2
    68 = 1
    69 = 0
    70 = 1
    71 = 1
    72 = 2
    73 = 1
    74 = 2
    76 = 1
    77 = 2
    78 = 1
    79 = 1
    80 = 1
    82 = 2
    83 = 2
    85 = 3
    86 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
